# 🤖 自动更新价格数据 2025-12-21 04:07:46
#
# The source data table keeps a rolling daily history of commodity prices
# in descending date order starting at row 2. A new trading day's data is
# published by inserting a fresh row right under the header (row 2),
# pushing all the previous rows down by one, and filling the new row with
# the latest date (one day after the prior newest date) together with
# that day's prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2 (above the previous newest entry),
# shifting all existing data rows down by one.
$ws.Rows("2:2").Insert()

# The freshly inserted row can inherit styling/number-formatting from
# neighboring rows (e.g. bold, or auto-detected date formats); clear it so
# the new row matches the plain/unstyled look of the rest of the data rows.
$ws.Range("A2:D2").ClearFormats()

# Force column A to be treated as literal text so the date string is not
# auto-converted into a numeric date serial by Excel.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-12-21"

# Remove the temporary text number-format override once the text value is
# safely stored, returning the cell to the default (unstyled) look.
$ws.Range("A2").ClearFormats()

$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
